$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Municipal Zone Summary")
$ws.Cells.Item(2, 4).Value = 155.3578315248243
$ws.Cells.Item(2, 7).Value = 260.5158083523038
$ws.Cells.Item(3, 4).Value = 279.4997003628848
$ws.Cells.Item(3, 7).Value = 785.6449465616154
$ws.Cells.Item(4, 4).Value = 46.47836900713188
$ws.Cells.Item(4, 7).Value = 91.87176664801527
$ws.Cells.Item(5, 4).Value = 27.12191921329827
$ws.Cells.Item(5, 7).Value = 42.34514907366468
$ws.Cells.Item(6, 4).Value = 90.32100206688119
$ws.Cells.Item(6, 7).Value = 181.0218630137317
$ws.Cells.Item(7, 4).Value = 88.83820760501696
$ws.Cells.Item(7, 7).Value = 156.0973323957411
$ws.Cells.Item(8, 4).Value = 124.0634955533343
$ws.Cells.Item(8, 7).Value = 330.4634125965265
$ws.Cells.Item(9, 4).Value = 15.53705587946277
$ws.Cells.Item(9, 7).Value = 20.32971258737049
$ws.Cells.Item(10, 4).Value = 61.06399030587847
$ws.Cells.Item(10, 7).Value = 81.53694608399567
$ws.Cells.Item(11, 4).Value = 17.23362346163252
$ws.Cells.Item(11, 7).Value = 29.29661505819229
$ws.Cells.Item(12, 4).Value = 50.11702191503642
$ws.Cells.Item(12, 7).Value = 186.6480070359511
$ws.Cells.Item(13, 4).Value = 137.8518712211186
$ws.Cells.Item(13, 7).Value = 231.6849390244221
$ws.Cells.Item(14, 4).Value = 56.1836652679048
$ws.Cells.Item(14, 7).Value = 93.91208427609222
$ws.Cells.Item(15, 4).Value = 320.6245116209929
$ws.Cells.Item(15, 7).Value = 742.7761373426529
$ws.Cells.Item(16, 4).Value = 2095.768818718121
$ws.Cells.Item(16, 7).Value = 4389.007663717412
$ws.Cells.Item(17, 4).Value = 75.49698962654575
$ws.Cells.Item(17, 7).Value = 190.9434917195673
$ws.Cells.Item(18, 4).Value = 20.35745769377773
$ws.Cells.Item(18, 7).Value = 34.79091803031537
$ws.Cells.Item(19, 4).Value = 13.82287693663128
$ws.Cells.Item(19, 7).Value = 21.86874613393865
$ws.Cells.Item(20, 4).Value = 37.36800264078069
$ws.Cells.Item(20, 7).Value = 114.7581598461426
$ws.Cells.Item(21, 4).Value = 373.6469898604245
$ws.Cells.Item(21, 7).Value = 1247.330569872149

$ws = $wb.Worksheets.Item("PMSA Summary")
$ws.Cells.Item(2, 4).Value = 173.83596285659
$ws.Cells.Item(2, 7).Value = 434.8040794200455
$ws.Cells.Item(3, 4).Value = 907.8183451033077
$ws.Cells.Item(3, 7).Value = 1858.356790161298
$ws.Cells.Item(4, 4).Value = 1176.797300265687
$ws.Cells.Item(4, 7).Value = 2513.78311302086
$ws.Cells.Item(5, 4).Value = 332.0095273340321
$ws.Cells.Item(5, 7).Value = 856.2437450898512
$ws.Cells.Item(6, 4).Value = 365.7409831235498
$ws.Cells.Item(6, 7).Value = 756.5429240343315
$ws.Cells.Item(7, 4).Value = 464.7854455783876
$ws.Cells.Item(7, 7).Value = 817.4292232033445
$ws.Cells.Item(8, 4).Value = 281.5214368060617
$ws.Cells.Item(8, 7).Value = 728.1156113856521
$ws.Cells.Item(9, 4).Value = 10.59740955363887
$ws.Cells.Item(9, 7).Value = 20.23821318226879
$ws.Cells.Item(10, 4).Value = 373.6469898604245
$ws.Cells.Item(10, 7).Value = 1247.330569872149
